$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# --- Fix existing rows ---

# Row 13 (ADDU / 0xFFFFFFFF + 0x1): outcome becomes an overflow error instead of 0
$ws.Range("D13").Value = "Error: Arithmetic overflow"

# Row 15 (ADDU / 0x7FFFFFFF + 0x7FFFFFFF): corrected outcome value
$ws.Range("D15").Value = -2

# Row 16 (ADDU / 0x8000000 + 0x1): corrected outcome value
$ws.Range("D16").Value = -2147483647

# Row 17 (ADDU / 0x1 + 0x80000000): corrected outcome value
$ws.Range("D17").Value = -2147483647

# Row 23 test case text changed from "0x8000000 + 0x80000000" to "0x8000000 + 0x8000"
$ws.Range("C23").Value = "0x8000000 + 0x8000"
$ws.Range("D23").Value = "Error: Arithmetic overflow"

# --- New rows 24-30: ADDIU (Add immediate unsigned, no overflow) test block ---

$ws.Range("A24").Value = 6
$ws.Range("B24").Value = "ADDIU"
$ws.Range("C24").Value = "0xFFFFFFFF + 0x1"
$ws.Range("D24").Value = "Error: Arithmetic overflow"

$ws.Range("C25").Value = "0x7000000 + 0xFFFF"
$ws.Range("D25").Value = 1879113727
$ws.Range("D7").Copy()
$ws.Range("D25").PasteSpecial(-4122)

$ws.Range("C26").Value = "0x7FFFFFFF + 0x7FFF"
$ws.Range("D26").Value = -2147450882

$ws.Range("C27").Value = "0x8000000 + 0x1"
$ws.Range("D27").Value = -2147483647

$ws.Range("C28").Value = "0x1 + 0x80000000"
$ws.Range("D28").Value = -2147483647

$ws.Range("C29").Value = "0xFFFFFFFF + 0x80000001"
$ws.Range("D29").Value = "Error: Arithmetic overflow"

$ws.Range("C30").Value = "0x8000000 + 0x80000000"
$ws.Range("D30").Value = "Error: Arithmetic overflow"

# --- View state: selection / scroll position ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E30").Select()
